$d = $word.ActiveDocument

# 1. "Uniform Resource API" line: append messaging detail to "Reactive dataflow."
$d.Content.Find.Execute(
    "Uniform Resource API: Sets, FCA, DOM layers, Monads. Reactive dataflow.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uniform Resource API: Sets, FCA, DOM layers, Monads. Reactive message driven dataflow (topics / signatures).",
    2)

# 2. "DOM Layers / OntResource hierarchy (FCA augmentation)." -> drop the
#    parenthetical qualifier.
$d.Content.Find.Execute(
    "DOM Layers / OntResource hierarchy (FCA augmentation).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "DOM Layers / OntResource hierarchy.",
    2)

# 3. Insert a brand-new paragraph (with its own blank spacer paragraph, matching
#    the surrounding blank-line-separated layout) right after that paragraph and
#    before the existing blank spacer that precedes "Parse DOM: ...".
$anchor = $d.Content
$anchor.Find.Execute("DOM Layers / OntResource hierarchy.")
$anchorPara = $anchor.Paragraphs.Item(1)
$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

$newParaIndex = $anchorPara.Index + 2
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.InsertBefore(
    "(Sets aggregation populates DOM layers FCA augmented or Sets aggregation builds FCA contexts rendered into FCA augmented DOM layers)."
)

# 4. "Parse DOM: ..." paragraph becomes a much longer "Functors. Parse DOM: ..."
#    paragraph.
$d.Content.Find.Execute(
    "Parse DOM: Relationship / Entity Monads (selectors / contexts)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Functors. Parse DOM: Instantiate Relationship / Entity Monads (selectors / contexts). Model services interactions renders functors possible transforms as browseable (HATEOAS reified) resources / contexts: reactive dialogs / prompts (HATEOAS / HAL protocols).",
    2)

# 5. "Model Services: ..." paragraph gets reworded.
$d.Content.Find.Execute(
    "Model Services: Monads parsed DOM interactions services (contexts). Render / update DOM.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Model Services: Browse DOM layers. Monads parsed DOM interactions services (functor contexts) available as operations over rendered models (HATEOAS).",
    2)

# 6. "Interactions: ..." paragraph: lower-case the second "Browse" and append a
#    new trailing sentence.
$d.Content.Find.Execute(
    "Interactions: Services. Browse DOM. Apply selectors / Browse available transforms (Monads / HATEOAS).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Interactions: Services. Browse DOM. Apply selectors / browse available transforms (Monads / HATEOAS). Monads applications render / update DOM / HATEOAS browsing response.",
    2)

# 7. "Outputs. ..." paragraph: "(events)" -> "(Events Inputs)".
$d.Content.Find.Execute(
    "Outputs. Connectors / Services (active Resource topics). Feedback (events).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Outputs. Connectors / Services (active Resource topics). Feedback (Events Inputs).",
    2)
